$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (Coin name) swaps ---
$ws.Range("B17").Value = "Chainlink"
$ws.Range("B18").Value = "TRON"
$ws.Range("B34").Value = "Hedera"
$ws.Range("B35").Value = "OKB"
$ws.Range("B40").Value = "TheGraph"
$ws.Range("B41").Value = "FirstDigitalUSD"

# --- Column C (Link) swaps ---
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"

# --- Column D (Price) updates. Leading apostrophe forces text entry;
#     Style reset back to Normal avoids leaving a Text-format style behind. ---
$ws.Range("D2").Value = "'67.490.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'3.667.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Value = "'579.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'169.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'3.656.82"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.619"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").Value = "'0.698"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.160"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'50.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.0000286"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'10.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'4.206.86"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'3.640.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'19.25"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.126"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'12.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'1.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'67.107.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'403.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'4.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'87.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'3.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'12.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'10.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = "'5.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Value = "'9.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Value = "'32.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Value = "'7.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = "'12.30"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "'0.115"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "'64.07"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = "'42.91"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = "'587.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Value = "'0.0₃0883"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Value = "'0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = "'0.391"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.998"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.133"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'2.95"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'2.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.0431"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'2.82"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'9.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'2.755.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.133"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'3.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'2.64"
$ws.Range("D51").Style = "Normal"

# --- Column E (Volume 1h %) updates ---
$ws.Range("E2").Value = "  -7.15%  "
$ws.Range("E3").Value = "  -7.10%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("E5").Value = "  -4.24%  "
$ws.Range("E6").Value = "  -0.84%  "
$ws.Range("E7").Value = "  -7.19%  "
$ws.Range("E8").Value = "  -9.16%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("E10").Value = "  -11.41%  "
$ws.Range("E11").Value = "  -11.47%  "
$ws.Range("E12").Value = "  -9.39%  "
$ws.Range("E13").Value = "  -12.46%  "
$ws.Range("E14").Value = "  -10.84%  "
$ws.Range("E15").Value = "  -8.12%  "
$ws.Range("E16").Value = "  -8.01%  "
$ws.Range("E17").Value = "  -10.11%  "
$ws.Range("E18").Value = "  -3.67%  "
$ws.Range("E19").Value = "  -9.39%  "
$ws.Range("E20").Value = "  -9.96%  "
$ws.Range("E21").Value = "  -7.65%  "
$ws.Range("E22").Value = "  -9.48%  "
$ws.Range("E23").Value = "  -8.49%  "
$ws.Range("E24").Value = "  -8.84%  "
$ws.Range("E25").Value = "  -9.83%  "
$ws.Range("E26").Value = "  -11.07%  "
$ws.Range("E27").Value = "  -6.15%  "
$ws.Range("E28").Value = "  +1.58%  "
$ws.Range("E29").Value = "  -12.04%  "
$ws.Range("E30").Value = "  -9.93%  "
$ws.Range("E31").Value = "  -10.00%  "
$ws.Range("E32").Value = "  -6.00%  "
$ws.Range("E33").Value = "  -11.50%  "
$ws.Range("E34").Value = "  -10.44%  "
$ws.Range("E35").Value = "  -7.39%  "
$ws.Range("E36").Value = "  -13.59%  "
$ws.Range("E37").Value = "  -7.34%  "
$ws.Range("E38").Value = "  -11.17%  "
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("E40").Value = "  -8.88%  "
$ws.Range("E41").Value = "  -0.42%  "
$ws.Range("E42").Value = "  -9.13%  "
$ws.Range("E43").Value = "  -14.12%  "
$ws.Range("E44").Value = "  +1.72%  "
$ws.Range("E45").Value = "  -10.02%  "
$ws.Range("E46").Value = "  -11.37%  "
$ws.Range("E47").Value = "  -14.56%  "
$ws.Range("E48").Value = "  -2.60%  "
$ws.Range("E49").Value = "  -10.38%  "
$ws.Range("E50").Value = "  -7.24%  "
$ws.Range("E51").Value = "  -6.99%  "
